$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 225 (weekly update adds two
# new price observations "ahead" of the existing ones), shifting the old
# rows 225-237 down to 227-239.
$ws.Rows.Item(225).Insert()
$ws.Rows.Item(225).Insert()

# --- New row 225 ---
$ws.Cells.Item(225, 1).Value = 6
$ws.Cells.Item(225, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(225, 3).Value = "Metropolitana"
$ws.Cells.Item(225, 4).Value = 44706
$ws.Cells.Item(225, 5).Value = 13
$ws.Cells.Item(225, 6).Value = 100112026
$ws.Cells.Item(225, 7).Value = "Haba"
$ws.Cells.Item(225, 8).Value = "Sin especificar"
$ws.Cells.Item(225, 9).Value = "Primera"
$ws.Cells.Item(225, 10).Value = 260
$ws.Cells.Item(225, 11).Value = 18000
$ws.Cells.Item(225, 12).Value = 19000
$ws.Cells.Item(225, 13).Value = 18462
$ws.Cells.Item(225, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(225, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(225, 16).Value = 738
$ws.Cells.Item(225, 17).Value = 25
$ws.Cells.Item(225, 18).Value = "Hortaliza"

# --- New row 226 ---
$ws.Cells.Item(226, 1).Value = 6
$ws.Cells.Item(226, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(226, 3).Value = "Metropolitana"
$ws.Cells.Item(226, 4).Value = 44706
$ws.Cells.Item(226, 5).Value = 13
$ws.Cells.Item(226, 6).Value = 100112026
$ws.Cells.Item(226, 7).Value = "Haba"
$ws.Cells.Item(226, 8).Value = "Sin especificar"
$ws.Cells.Item(226, 9).Value = "Primera"
$ws.Cells.Item(226, 10).Value = 180
$ws.Cells.Item(226, 11).Value = 16000
$ws.Cells.Item(226, 12).Value = 17000
$ws.Cells.Item(226, 13).Value = 16444
$ws.Cells.Item(226, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(226, 15).Value = "Región Metropolitana"
$ws.Cells.Item(226, 16).Value = 658
$ws.Cells.Item(226, 17).Value = 25
$ws.Cells.Item(226, 18).Value = "Hortaliza"
